$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift the START DATE / END DATE values forward by the same offset
# (0.548086442395288 days, ~13h09m14s) for rows 2 and 3.
$ws.Range("F2").Value = 45385.89280657649
$ws.Range("G2").Value = 45387.89280657649
$ws.Range("F3").Value = 45387.89280657649
$ws.Range("G3").Value = 45389.89280657649
